# Auto-generated edit script applying the cryptos.xlsx symbol-list refresh
# (commit: "Updated symbol list on Thu Jan 12 23:56:59 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each assignment is prefixed with a leading apostrophe so Excel stores the
# value as literal text (matching the original t="inlineStr" cells) instead
# of auto-converting number- or percent-looking strings into numeric values.

$ws.Range('D2').Value = "'287.56"
$ws.Range('E2').Value = "'1.43%"
$ws.Range('D3').Value = "'29.61"
$ws.Range('E3').Value = "'7.32%"
$ws.Range('D4').Value = "'5.082"
$ws.Range('E4').Value = "'4.27%"
$ws.Range('D5').Value = "'0.06658"
$ws.Range('E5').Value = "'2.17%"
$ws.Range('D6').Value = "'7.393"
$ws.Range('E6').Value = "'3.87%"
$ws.Range('B7').Value = "'FTXToken"
$ws.Range('C7').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D7').Value = "'1.378"
$ws.Range('E7').Value = "'8.12%"
$ws.Range('B8').Value = "'MXToken"
$ws.Range('C8').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D8').Value = "'0.9382"
$ws.Range('E8').Value = "'2.65%"
$ws.Range('B9').Value = "'WazirX"
$ws.Range('C9').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D9').Value = "'0.1579"
$ws.Range('E9').Value = "'1.93%"
$ws.Range('B10').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C10').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D10').Value = "'0.06538"
$ws.Range('E10').Value = "'3.26%"
$ws.Range('B11').Value = "'MandalaExchangeToken"
$ws.Range('C11').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D11').Value = "'0.07629"
$ws.Range('E11').Value = "'1.43%"
$ws.Range('B12').Value = "'BitrueCoin"
$ws.Range('C12').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D12').Value = "'0.02935"
$ws.Range('E12').Value = "'0.77%"
$ws.Range('B13').Value = "'BitMartToken"
$ws.Range('C13').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D13').Value = "'0.08988"
$ws.Range('E13').Value = "'0.27%"
$ws.Range('B14').Value = "'BitForexToken"
$ws.Range('C14').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D14').Value = "'0.001598"
$ws.Range('E14').Value = "'0.15%"
$ws.Range('B15').Value = "'CoinExToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D15').Value = "'0.04496"
$ws.Range('E15').Value = "'2.36%"
$ws.Range('B16').Value = "'One"
$ws.Range('C16').Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('D16').Value = "'0.0006482"
$ws.Range('E16').Value = "'0.44%"
$ws.Range('B17').Value = "'TigerCash"
$ws.Range('C17').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('D17').Value = "'0.006296"
$ws.Range('E17').Value = "'4.85%"
$ws.Range('B18').Value = "'LEO"
$ws.Range('C18').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('D18').Value = "'3.446"
$ws.Range('E18').Value = "'-1.74%"
$ws.Range('B19').Value = "'GateToken"
$ws.Range('C19').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D19').Value = "'3.404"
$ws.Range('E19').Value = "'2.07%"
$ws.Range('E20').Value = "'1.08%"
$ws.Range('D21').Value = "'0.3215"
$ws.Range('E21').Value = "'2.21%"
$ws.Range('E22').Value = "'-3.78%"
$ws.Range('D23').Value = "'4.091"
$ws.Range('E23').Value = "'5.30%"
$ws.Range('D24').Value = "'0.1553"
$ws.Range('E24').Value = "'3.82%"
$ws.Range('E25').Value = "'1.25%"
$ws.Range('D26').Value = "'0.004135"
$ws.Range('E26').Value = "'-4.05%"
$ws.Range('D27').Value = "'0.0001250"
$ws.Range('E27').Value = "'6.59%"
$ws.Range('D28').Value = "'0.0001618"
$ws.Range('E28').Value = "'-1.88%"
$ws.Range('D40').Value = "'0.04220"
$ws.Range('E40').Value = "'2.43%"
$ws.Range('D41').Value = "'0.006716"
$ws.Range('E41').Value = "'0.98%"
$ws.Range('D42').Value = "'0.1247"
$ws.Range('E42').Value = "'-10.01%"
$ws.Range('D43').Value = "'0.001969"
$ws.Range('E43').Value = "'-5.16%"
$ws.Range('D44').Value = "'0.01240"
$ws.Range('E44').Value = "'6.29%"
$ws.Range('D45').Value = "'0.00005568"
$ws.Range('E45').Value = "'1.10%"
$ws.Range('E46').Value = "'20.74%"
$ws.Range('D47').Value = "'0.01307"
$ws.Range('E47').Value = "'-28.96%"
